$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.014.02"
$ws.Range("E2").Value = "  -0.17%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.95"
$ws.Range("E3").Value = "  -0.43%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.77"
$ws.Range("E5").Value = "  -0.93%  "

# Row 6
$ws.Range("E6").Value = "  -0.74%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  -1.75%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0623"
$ws.Range("E9").Value = "  -2.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.34"
$ws.Range("E10").Value = "  -6.79%  "

# Row 11
$ws.Range("E11").Value = "  -0.83%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.866.69"
$ws.Range("E12").Value = "  -0.37%  "

# Row 13
$ws.Range("E13").Value = "  -2.06%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.632.46"
$ws.Range("E14").Value = "  -0.61%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.528"
$ws.Range("E15").Value = "  -3.06%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.019.88"
$ws.Range("E16").Value = "  -0.47%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0744"
$ws.Range("E17").Value = "  -2.82%  "

# Row 18
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.84"
$ws.Range("E18").Value = "  -2.48%  "

# Row 19
$ws.Range("E19").Value = "  -0.10%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.99"
$ws.Range("E20").Value = "  -0.77%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.24"
$ws.Range("E21").Value = "  -2.33%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.73"
$ws.Range("E22").Value = "  -2.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.10"
$ws.Range("E23").Value = "  -2.25%  "

# Row 24
$ws.Range("E24").Value = "  +2.28%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.96"
$ws.Range("E25").Value = "  -0.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.79"
$ws.Range("E26").Value = "  -1.30%  "

# Row 27
$ws.Range("E27").Value = "  -0.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.83"
$ws.Range("E28").Value = "  -1.12%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.23"
$ws.Range("E29").Value = "  -2.10%  "

# Row 30
$ws.Range("E30").Value = "  -1.33%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0484"
$ws.Range("E31").Value = "  -2.67%  "

# Row 32
$ws.Range("E32").Value = "  -3.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.14"
$ws.Range("E33").Value = "  -4.34%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.41"
$ws.Range("E34").Value = "  -1.98%  "

# Row 35
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  -2.42%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.136.80"
$ws.Range("E36").Value = "  +0.23%  "

# Row 37
$ws.Range("E37").Value = "  -4.09%  "

# Row 38
$ws.Range("E38").Value = "  -0.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.520"
$ws.Range("E39").Value = "  -3.79%  "

# Row 40
$ws.Range("E40").Value = "  -1.32%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.50"
$ws.Range("E41").Value = "  -1.14%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.780"
$ws.Range("E42").Value = "  -2.34%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.776.19"
$ws.Range("E43").Value = "  -0.44%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.23"
$ws.Range("E44").Value = "  -4.90%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0116"
$ws.Range("E45").Value = "  -1.36%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.12"
$ws.Range("E46").Value = "  -2.85%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0528"
$ws.Range("E47").Value = "  -0.34%  "

# Row 48
$ws.Range("E48").Value = "  +2.27%  "

# Row 49
$ws.Range("E49").Value = "  -0.42%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.55"
$ws.Range("E50").Value = "  -2.59%  "

# Row 51
$ws.Range("E51").Value = "  +0.00%  "
